$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.456.52'
$ws.Range('E2').Value = '  +2.32%  '
$ws.Range('D3').Value = '3.077.55'
$ws.Range('E3').Value = '  +4.35%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''577.97'
$ws.Range('E5').Value = '  +1.54%  '
$ws.Range('D6').Value = '''167.33'
$ws.Range('E6').Value = '  +5.28%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '3.076.05'
$ws.Range('E8').Value = '  +4.47%  '
$ws.Range('D9').Value = '''0.523'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').Value = '''0.481'
$ws.Range('E12').Value = '  +5.57%  '
$ws.Range('D13').Value = '''0.0000248'
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('D14').Value = '''36.36'
$ws.Range('E14').Value = '  +6.60%  '
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').Value = '3.588.04'
$ws.Range('E16').Value = '  +4.32%  '
$ws.Range('D17').Value = '66.590.58'
$ws.Range('E17').Value = '  +2.39%  '
$ws.Range('D18').Value = '''7.21'
$ws.Range('E18').Value = '  +4.04%  '
$ws.Range('D19').Value = '3.078.20'
$ws.Range('E19').Value = '  +4.34%  '
$ws.Range('D20').Value = '''16.13'
$ws.Range('E20').Value = '  +16.55%  '
$ws.Range('D21').Value = '''466.92'
$ws.Range('E21').Value = '  +4.69%  '
$ws.Range('D22').Value = '''0.713'
$ws.Range('E22').Value = '  +4.70%  '
$ws.Range('D23').Value = '''7.52'
$ws.Range('E23').Value = '  +3.81%  '
$ws.Range('D24').Value = '''83.04'
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('E25').Value = '  +4.80%  '
$ws.Range('D26').Value = '''12.85'
$ws.Range('E26').Value = '  +6.97%  '
$ws.Range('D27').Value = '''10.10'
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D29').Value = '''7.96'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('E31').Value = '  +2.88%  '
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('D33').Value = '''28.15'
$ws.Range('E33').Value = '  +3.58%  '
$ws.Range('E34').Value = '  +4.05%  '
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  +2.46%  '
$ws.Range('D37').Value = '''5.87'
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('E38').Value = '  +6.91%  '
$ws.Range('D39').Value = '''46.26'
$ws.Range('E39').Value = '  +5.38%  '
$ws.Range('D40').Value = '''50.17'
$ws.Range('E40').Value = '  +2.57%  '
$ws.Range('D41').Value = '''0.317'
$ws.Range('E41').Value = '  +6.17%  '
$ws.Range('E42').Value = '  +2.54%  '
$ws.Range('D43').Value = '''8.66'
$ws.Range('E43').Value = '  +3.14%  '
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('D46').Value = '''381.87'
$ws.Range('E46').Value = '  -0.74%  '
$ws.Range('D47').Value = '2.758.69'
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('D48').Value = '''134.81'
$ws.Range('E48').Value = '  +1.95%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '''24.58'
$ws.Range('E50').Value = '  +5.87%  '
$ws.Range('D51').Value = '''2.23'
$ws.Range('E51').Value = '  +3.84%  '
